# Resolve bug with recursion: split the single "Seed Device:" prompt on the
# Audit sheet into two separate seed-device prompts, "Seed Device 1:" and
# "Seed Device 2:".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Audit")

# Rename existing "Seed Device:" label (A7) to "Seed Device 1:"
$ws.Range("A7").Value = "Seed Device 1:"

# Add a new "Seed Device 2:" label in the previously blank cell below (A8)
$ws.Range("A8").Value = "Seed Device 2:"

# Reflect the user's last active cell being A8 after the edit
$ws.Activate()
$ws.Range("A8").Select()
